# Adding labs 11 and 13
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 ---
$ws.Range("B4").Value = 97
$ws.Range("C4").Value = 87

# --- Row 5 ---
$ws.Range("B5").Value = 5346
$ws.Range("C5").Value = 4777
$ws.Range("D5").Value = 4597
$ws.Range("E5").Value = 4607
$ws.Range("I5").Value = 125
$ws.Range("J5").Value = 136
$ws.Range("K5").Value = 146
$ws.Range("L5").Value = 182

# --- Row 6 ---
$ws.Range("C6").Value = 74
$ws.Range("D6").Value = 70

# --- Row 7 ---
$ws.Range("B7").Value = 75
$ws.Range("D7").Value = 65
$ws.Range("E7").Value = 66
$ws.Range("H7").Value = 0.1673611111111111

# --- Row 8 ---
$ws.Range("C8").Value = 154
$ws.Range("D8").Value = 153
$ws.Range("E8").Value = 144
$ws.Range("H8").Value = 2120

# --- Row 9 ---
$ws.Range("B9").Value = 3.4
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 7.9
$ws.Range("E9").Value = 6.8

# --- Row 10 ---
$ws.Range("C10").Value = 0.42
$ws.Range("E10").Value = 0.46

# --- Row 12 ---
$ws.Range("D12").Value = 2.5
$ws.Range("E12").Value = 2.2000000000000002

# --- Row 14 ---
$ws.Range("C14").Value = 0.09
$ws.Range("D14").Value = 0.04
$ws.Range("E14").Value = 0.11

# --- Row 16 ---
$ws.Range("B16").Value = 21
$ws.Range("C16").Value = 20
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 8

# --- Row 17 ---
$ws.Range("B17").Value = 5413
$ws.Range("C17").Value = 5409
$ws.Range("D17").Value = 5408
$ws.Range("E17").Value = 5396

# --- Row 18 ---
$ws.Range("B18").Value = 2451
$ws.Range("C18").Value = 2451
$ws.Range("D18").Value = 2451
$ws.Range("E18").Value = 2421

# --- Row 19 ---
$ws.Range("B19").Value = 2962
$ws.Range("C19").Value = 2958
$ws.Range("D19").Value = 2958
$ws.Range("E19").Value = 2976

# --- Row 20 ---
$ws.Range("B20").Value = 45
$ws.Range("C20").Value = 45
$ws.Range("D20").Value = 45
$ws.Range("E20").Value = 45

# --- View: scroll/selection change ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A10").Select()
